# Auto-generated edit script applying the cryptos.xlsx diff (GitHub Actions crypto price refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.085.14"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.290.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.30"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.315.43"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0982"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.13"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.19%  "
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.78"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.700.52"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "55.125.17"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.286.41"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.47"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.11"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.65"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.154"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.27"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.38%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0712"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.63"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.08"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.23"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.919"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("E38").Value = "  +5.62%  "
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.377"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "136.31"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.12"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.53%  "
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "262.49"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +10.03%  "
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0915"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.555"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.375"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("E50").Value = "  +3.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.21%  "
